$wb = $excel.ActiveWorkbook

# Sheet "展览" (overview sheet 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 126
$ws1.Range("F3").Value = 2144
$ws1.Range("F4").Value = 25
$ws1.Range("F5").Value = 11215
$ws1.Range("F8").Value = 312
$ws1.Range("F10").Value = 11122
$ws1.Range("F11").Value = 448
$ws1.Range("F13").Value = 48
$ws1.Range("F14").Value = 1727
$ws1.Range("F15").Value = 5572
$ws1.Range("F17").Value = 3446
$ws1.Range("F18").Value = 7

# Sheet "全部类型" (combined overview sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 126
$ws4.Range("F3").Value = 2144
$ws4.Range("F5").Value = 25
$ws4.Range("F7").Value = 11215
$ws4.Range("F10").Value = 312
$ws4.Range("F12").Value = 11122
$ws4.Range("F13").Value = 448
$ws4.Range("F15").Value = 48
$ws4.Range("F16").Value = 1727
$ws4.Range("F17").Value = 5572
$ws4.Range("F19").Value = 3446
$ws4.Range("F20").Value = 7
